$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 10798
$ws1.Range("F4").Value = 67
$ws1.Range("F5").Value = 725
$ws1.Range("F6").Value = 500

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 6

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 10798
$ws4.Range("F4").Value = 67
$ws4.Range("F5").Value = 725
$ws4.Range("F6").Value = 6
$ws4.Range("F7").Value = 500

$wb.Save()
